$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.29094123840332
$ws.Range("B1").Value = 6.618060111999512
$ws.Range("C1").Value = 6.490407943725586
$ws.Range("D1").Value = 6.864789009094238
$ws.Range("E1").Value = 3.450332641601562
